$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q6").Value = 25.4
$ws.Range("S6").Value = 33.2
$ws.Range("S6").Copy()
$ws.Range("T6").PasteSpecial(-4122)
$ws.Range("T6").Value = 36
$ws.Range("S7").Copy()
$ws.Range("T7").PasteSpecial(-4122)
$ws.Range("T7").Value = 51.4
$ws.Range("S8").Copy()
$ws.Range("T8").PasteSpecial(-4122)
$ws.Range("T8").Value = 15.1
$ws.Range("R12").Value = 0.46
$ws.Range("R12").Copy()
$ws.Range("S12").PasteSpecial(-4122)
$ws.Range("S12").Value = 0.82
$ws.Range("R19").Copy()
$ws.Range("S19").PasteSpecial(-4122)
$ws.Range("S19").Value = 13.83
$ws.Range("R19").Copy()
$ws.Range("T19").PasteSpecial(-4122)
$ws.Range("T19").Value = 14.74
$ws.Range("R26").Copy()
$ws.Range("S26").PasteSpecial(-4122)
$ws.Range("S26").Value = 426.2
$ws.Range("R27").Copy()
$ws.Range("S27").PasteSpecial(-4122)
$ws.Range("S27").Value = 406.2
$ws.Range("R28").Copy()
$ws.Range("S28").PasteSpecial(-4122)
$ws.Range("S28").Value = 444.9
$ws.Range("R29").Copy()
$ws.Range("S29").PasteSpecial(-4122)
$ws.Range("S29").Value = 430.6
$ws.Range("R30").Copy()
$ws.Range("S30").PasteSpecial(-4122)
$ws.Range("S30").Value = 419.7
$ws.Range("R31").Copy()
$ws.Range("S31").PasteSpecial(-4122)
$ws.Range("S31").Value = 253.9
$ws.Range("R32").Copy()
$ws.Range("S32").PasteSpecial(-4122)
$ws.Range("S32").Value = 283.1
$ws.Range("R33").Copy()
$ws.Range("S33").PasteSpecial(-4122)
$ws.Range("S33").Value = 226.6
$ws.Range("R34").Copy()
$ws.Range("S34").PasteSpecial(-4122)
$ws.Range("S34").Value = 302
$ws.Range("R35").Copy()
$ws.Range("S35").PasteSpecial(-4122)
$ws.Range("S35").Value = 222.2
$ws.Range("R36").Copy()
$ws.Range("S36").PasteSpecial(-4122)
$ws.Range("S36").Value = 28.8
$ws.Range("R37").Copy()
$ws.Range("S37").PasteSpecial(-4122)
$ws.Range("S37").Value = 27.2
$ws.Range("R38").Copy()
$ws.Range("S38").PasteSpecial(-4122)
$ws.Range("S38").Value = 30.3
$ws.Range("R39").Copy()
$ws.Range("S39").PasteSpecial(-4122)
$ws.Range("S39").Value = 29.2
$ws.Range("R40").Copy()
$ws.Range("S40").PasteSpecial(-4122)
$ws.Range("S40").Value = 28.2
$ws.Range("R41").Copy()
$ws.Range("S41").PasteSpecial(-4122)
$ws.Range("S41").Value = 26.6
$ws.Range("R42").Copy()
$ws.Range("S42").PasteSpecial(-4122)
$ws.Range("S42").Value = 31.3
$ws.Range("R43").Copy()
$ws.Range("S43").PasteSpecial(-4122)
$ws.Range("S43").Value = 22.2
$ws.Range("R44").Copy()
$ws.Range("S44").PasteSpecial(-4122)
$ws.Range("S44").Value = 27.9
$ws.Range("R45").Copy()
$ws.Range("S45").PasteSpecial(-4122)
$ws.Range("S45").Value = 24.6
$ws.Range("H6").Copy()
$ws.Range("K65").PasteSpecial(-4122)
$ws.Range("K65").Value = 23.3
$ws.Range("H6").Copy()
$ws.Range("M65").PasteSpecial(-4122)
$ws.Range("M65").Value = 35.3
$ws.Range("H6").Copy()
$ws.Range("O65").PasteSpecial(-4122)
$ws.Range("O65").Value = 56.6
$ws.Range("H6").Copy()
$ws.Range("R65").PasteSpecial(-4122)
$ws.Range("R65").Value = 68.4
$ws.Range("H6").Copy()
$ws.Range("S65").PasteSpecial(-4122)
$ws.Range("S65").Value = 89
$ws.Range("H6").Copy()
$ws.Range("T65").PasteSpecial(-4122)
$ws.Range("T65").Value = 83.3
$ws.Range("Q118").Copy()
$ws.Range("S118").PasteSpecial(-4122)
$ws.Range("S118").Value = 46.1
$ws.Range("Q119").Copy()
$ws.Range("S119").PasteSpecial(-4122)
$ws.Range("S119").Value = 33.1
$ws.Range("Q120").Copy()
$ws.Range("S120").PasteSpecial(-4122)
$ws.Range("S120").Value = 42.4
$ws.Range("Q121").Copy()
$ws.Range("S121").PasteSpecial(-4122)
$ws.Range("S121").Value = 32.3
$ws.Range("Q122").Copy()
$ws.Range("S122").PasteSpecial(-4122)
$ws.Range("S122").Value = 56.3
$ws.Range("Q123").Copy()
$ws.Range("S123").PasteSpecial(-4122)
$ws.Range("S123").Value = 35.4
$ws.Range("Q124").Copy()
$ws.Range("S124").PasteSpecial(-4122)
$ws.Range("S124").Value = 59.7
$ws.Range("Q125").Value = 25.4
$ws.Range("S125").Value = 33.2
$ws.Range("S125").Copy()
$ws.Range("T125").PasteSpecial(-4122)
$ws.Range("T125").Value = 36
$ws.Range("S126").Copy()
$ws.Range("T126").PasteSpecial(-4122)
$ws.Range("T126").Value = 51.4
$ws.Range("S127").Copy()
$ws.Range("T127").PasteSpecial(-4122)
$ws.Range("T127").Value = 15.1
$ws.Range("R138").Value = 15.61
$ws.Range("R138").Copy()
$ws.Range("S138").PasteSpecial(-4122)
$ws.Range("S138").Value = 16.88
$ws.Range("R141").Copy()
$ws.Range("S141").PasteSpecial(-4122)
$ws.Range("S141").Value = 9.1
$ws.Range("R187").Copy()
$ws.Range("S187").PasteSpecial(-4122)
$ws.Range("S187").Value = 44.2
$ws.Range("S188").Copy()
$ws.Range("T188").PasteSpecial(-4122)
$ws.Range("T188").Value = 15.8
$ws.Range("S204").Copy()
$ws.Range("T204").PasteSpecial(-4122)
$ws.Range("T204").Value = 170
$ws.Range("P205").Copy()
$ws.Range("T205").PasteSpecial(-4122)
$ws.Range("T205").Value = 9.3
$ws.Range("Q207").Copy()
$ws.Range("T207").PasteSpecial(-4122)
$ws.Range("T207").Value = 6.1
$ws.Range("R216").Copy()
$ws.Range("S216").PasteSpecial(-4122)
$ws.Range("S216").Value = 14.8
$ws.Range("R222").Copy()
$ws.Range("S222").PasteSpecial(-4122)
$ws.Range("S222").Value = 7.6
$ws.Range("R226").Copy()
$ws.Range("S226").PasteSpecial(-4122)
$ws.Range("S226").Value = 64070
$ws.Range("N235").Value = 100.7
$ws.Range("R235").Copy()
$ws.Range("S235").PasteSpecial(-4122)
$ws.Range("S235").Value = 94.2
$ws.Range("L236").Value = 94.2
$ws.Range("R236").Copy()
$ws.Range("S236").PasteSpecial(-4122)
$ws.Range("S236").Value = 93.4
$ws.Range("R237").Value = 15.61
$ws.Range("R237").Copy()
$ws.Range("S237").PasteSpecial(-4122)
$ws.Range("S237").Value = 16.88
$ws.Range("Q238").Copy()
$ws.Range("R238").PasteSpecial(-4122)
$ws.Range("R238").Value = 1189
$ws.Range("Q238").Copy()
$ws.Range("S238").PasteSpecial(-4122)
$ws.Range("S238").Value = 1318
$ws.Range("S239").Copy()
$ws.Range("T239").PasteSpecial(-4122)
$ws.Range("T239").Value = 28779
$ws.Range("S247").Copy()
$ws.Range("T247").PasteSpecial(-4122)
$ws.Range("T247").Value = 75
$ws.Range("S248").Copy()
$ws.Range("T248").PasteSpecial(-4122)
$ws.Range("T248").Value = 7
$ws.Range("R249").Copy()
$ws.Range("S249").PasteSpecial(-4122)
$ws.Range("S249").Value = 80
$ws.Range("R249").Copy()
$ws.Range("T249").PasteSpecial(-4122)
$ws.Range("T249").Value = 97
$ws.Range("S289").Copy()
$ws.Range("T289").PasteSpecial(-4122)
$ws.Range("T289").Value = 2603.83
$ws.Range("S290").Copy()
$ws.Range("T290").PasteSpecial(-4122)
$ws.Range("T290").Value = 1744.48
$ws.Range("Q294").Copy()
$ws.Range("T294").PasteSpecial(-4122)
$ws.Range("T294").Value = 8
$ws.Range("S295").Copy()
$ws.Range("T295").PasteSpecial(-4122)
$ws.Range("T295").Value = 0
$ws.Range("S296").Copy()
$ws.Range("T296").PasteSpecial(-4122)
$ws.Range("T296").Value = 60
$excel.CutCopyMode = $false
